# Generate Report for Handoff
# Updates the localization-status report with a new handoff run:
#  - new source file GUID (65f972a3-... -> 4f1bf542-...)
#  - new xliff content hash (a1841e1e... -> ea97dbac1...)
#  - refreshed handoff/generate timestamps

$wb = $excel.ActiveWorkbook

$newGuid = "4f1bf542-ec5f-4c49-b5a4-5b0acbd09d0b"
$newHash = "ea97dbac1d4cda00b741e612c97eeda22bc29f0f"

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = "e2e\$newGuid.md"
foreach ($h in $ws.Hyperlinks) { $h.TextToDisplay = "e2e\$newGuid.md" }
$ws.Range("G2").Value = "2016-08-30 00:59:08"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "$newGuid.md"
foreach ($h in $ws.Hyperlinks) { $h.TextToDisplay = "$newGuid.md" }
$ws.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-30 00:58:59"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "$newGuid.md"
foreach ($h in $ws.Hyperlinks) { $h.TextToDisplay = "$newGuid.md" }
$ws.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$ws.Range("H2").Value = "2016-08-30 00:59:08"
